# The post previously at row 607 ("「自分を頼りなさい。宇宙で最も確かな物が
# 自分だと思って。」...") was removed from the source data. Deleting the
# entire worksheet row shifts every subsequent row up by one (rows 608-792
# become 607-791) and keeps all cell content, simply renumbering the rows,
# which is exactly the behavior captured in the diff. The sheet's dimension
# (A1:C792 -> A1:C791) is updated automatically by Excel when the row is
# removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(607).Delete()
